$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1557.8108
$ws.Range("I98").Value = 1169.6857
$ws.Range("J98").Value = 8350
$ws.Range("K98").Value = 1169.6857
$ws.Range("L98").Value = 8350
$ws.Range("M98").Value = 328.3143
$ws.Range("N98").Value = -11346

$ws.Range("H122").Value = 1557.8108
$ws.Range("I122").Value = 1169.6857
$ws.Range("J122").Value = 8350
$ws.Range("K122").Value = 3509.0571
$ws.Range("L122").Value = 25050
$ws.Range("M122").Value = -1059.0571
$ws.Range("N122").Value = -29950

$ws.Range("H132").Value = 4795.532
$ws.Range("I132").Value = 4005.6553
$ws.Range("K132").Value = 12016.9659
$ws.Range("M132").Value = -9486.965899999999

$ws.Range("H138").Value = 3409.5117
$ws.Range("I138").Value = 3877.8235
$ws.Range("K138").Value = 11633.4705
$ws.Range("M138").Value = -6493.470499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1910.8334
$ws.Range("I2").Value = 2010.1177
$ws.Range("K2").Value = 2010.1177
$ws.Range("M2").Value = -1897.1177

$ws.Range("H97").Value = 584.871
$ws.Range("I97").Value = 491.17648
$ws.Range("K97").Value = 491.17648
$ws.Range("M97").Value = 4.823519999999974

$ws.Range("H102").Value = 2626.8125
$ws.Range("I102").Value = 2626.8125
$ws.Range("K102").Value = 2626.8125
$ws.Range("M102").Value = -1004.8125

$ws.Range("H110").Value = 1286.0435
$ws.Range("I110").Value = 929
$ws.Range("K110").Value = 929
$ws.Range("M110").Value = 1116

$ws.Range("H116").Value = 1910.8334
$ws.Range("I116").Value = 2010.1177
$ws.Range("K116").Value = 2010.1177
$ws.Range("M116").Value = 283.8823

$ws.Range("H122").Value = 6278.357
$ws.Range("I122").Value = 5727.091
$ws.Range("K122").Value = 17181.273
$ws.Range("M122").Value = -14731.273

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 13676.714
$ws.Range("I132").Value = 15679.473
$ws.Range("J132").Value = 1660.1666
$ws.Range("K132").Value = 47038.419
$ws.Range("L132").Value = 4980.4998
$ws.Range("M132").Value = -44508.419
$ws.Range("N132").Value = -10040.4998

$ws.Range("H135").Value = 11448.75
$ws.Range("J135").Value = 11448.75
$ws.Range("L135").Value = 11448.75
$ws.Range("N135").Value = -21588.75

$ws.Range("H139").Value = 79949
$ws.Range("J139").Value = 79949
$ws.Range("L139").Value = 79949
$ws.Range("N139").Value = -90229

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1910.8334
$ws.Range("I3").Value = 2010.1177
$ws.Range("K3").Value = 2010.1177
$ws.Range("M3").Value = -1896.1177

$ws.Range("H80").Value = 898.5
$ws.Range("J80").Value = 898.5
$ws.Range("L80").Value = 898.5
$ws.Range("N80").Value = -2894.5

$ws.Range("H83").Value = 898.5
$ws.Range("J83").Value = 898.5
$ws.Range("L83").Value = 4492.5
$ws.Range("N83").Value = -14476.5

$ws.Range("H86").Value = 2831.4285
$ws.Range("I86").Value = 1626.68
$ws.Range("K86").Value = 1626.68
$ws.Range("M86").Value = -503.6800000000001

$ws.Range("H89").Value = 2831.4285
$ws.Range("I89").Value = 1626.68
$ws.Range("K89").Value = 8133.400000000001
$ws.Range("M89").Value = -2517.400000000001

$ws.Range("H94").Value = 1593.697
$ws.Range("I94").Value = 1471.25
$ws.Range("K94").Value = 1471.25
$ws.Range("M94").Value = -1020.25

$ws.Range("H105").Value = 3703.6667
$ws.Range("I105").Value = 3296.3333
$ws.Range("K105").Value = 3296.3333
$ws.Range("M105").Value = -1549.3333

$ws.Range("H107").Value = 2096.4285
$ws.Range("I107").Value = 1085.5
$ws.Range("K107").Value = 1085.5
$ws.Range("M107").Value = 834.5

$ws.Range("H134").Value = 2792.3235
$ws.Range("I134").Value = 2620.3635
$ws.Range("J134").Value = 3107.5833
$ws.Range("K134").Value = 7861.0905
$ws.Range("L134").Value = 9322.749899999999
$ws.Range("M134").Value = -5326.0905
$ws.Range("N134").Value = -14392.7499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1504.6666
$ws.Range("I16").Value = 1525.4
$ws.Range("K16").Value = 1525.4
$ws.Range("M16").Value = -1238.4

$ws.Range("H31").Value = 5005.1577
$ws.Range("I31").Value = 2984.375
$ws.Range("J31").Value = 6474.8184
$ws.Range("K31").Value = 2984.375
$ws.Range("L31").Value = 6474.8184
$ws.Range("M31").Value = -2689.375
$ws.Range("N31").Value = -7064.8184

$ws.Range("H34").Value = 5005.1577
$ws.Range("I34").Value = 2984.375
$ws.Range("J34").Value = 6474.8184
$ws.Range("K34").Value = 2984.375
$ws.Range("L34").Value = 6474.8184
$ws.Range("M34").Value = -2782.375
$ws.Range("N34").Value = -6878.8184

$ws.Range("H57").Value = 49999
$ws.Range("J57").Value = 49999
$ws.Range("L57").Value = 49999
$ws.Range("N57").Value = -51119

$ws.Range("H94").Value = 874.5
$ws.Range("I94").Value = 486.66666
$ws.Range("K94").Value = 486.66666
$ws.Range("M94").Value = -35.66665999999998

$ws.Range("H113").Value = 1504.6666
$ws.Range("I113").Value = 1525.4
$ws.Range("K113").Value = 1525.4
$ws.Range("M113").Value = 644.5999999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 1909.3334
$ws.Range("I104").Value = 1909.3334
$ws.Range("K104").Value = 5728.0002
$ws.Range("M104").Value = -3107.0002

$ws.Range("H131").Value = 13978.667
$ws.Range("J131").Value = 19217.445
$ws.Range("L131").Value = 57652.335
$ws.Range("N131").Value = -67732.33499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5833333.5
$ws.Range("J7").Value = 5833333.5
$ws.Range("L7").Value = 5833333.5
$ws.Range("N7").Value = -5833557.5

$ws.Range("H8").Value = 5833333.5
$ws.Range("J8").Value = 5833333.5
$ws.Range("L8").Value = 5833333.5
$ws.Range("N8").Value = -5833611.5

$ws.Range("H11").Value = 2316190.5
$ws.Range("J11").Value = 2316190.5
$ws.Range("L11").Value = 2316190.5
$ws.Range("N11").Value = -2316468.5

$ws.Range("H97").Value = 1196.9474
$ws.Range("I97").Value = 1328.3334
$ws.Range("J97").Value = 1078.7
$ws.Range("K97").Value = 1328.3334
$ws.Range("L97").Value = 1078.7
$ws.Range("M97").Value = -832.3334
$ws.Range("N97").Value = -2070.7

$ws.Range("H113").Value = 151239.28
$ws.Range("I113").Value = 133418.38
$ws.Range("K113").Value = 133418.38
$ws.Range("M113").Value = -131248.38

$ws.Range("H122").Value = 2366.4092
$ws.Range("I122").Value = 1978.35
$ws.Range("K122").Value = 5935.049999999999
$ws.Range("M122").Value = -3485.049999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4160.5
$ws.Range("I61").Value = 3540.5715
$ws.Range("K61").Value = 3540.5715
$ws.Range("M61").Value = -3338.5715

$ws.Range("H113").Value = 4160.5
$ws.Range("I113").Value = 3540.5715
$ws.Range("K113").Value = 3540.5715
$ws.Range("M113").Value = -1370.5715

$ws.Range("H132").Value = 36926.4
$ws.Range("I132").Value = 43325.414
$ws.Range("J132").Value = 5997.8335
$ws.Range("K132").Value = 129976.242
$ws.Range("L132").Value = 17993.5005
$ws.Range("M132").Value = -127446.242
$ws.Range("N132").Value = -23053.5005

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1126.3
$ws.Range("I100").Value = 783
$ws.Range("J100").Value = 2499.5
$ws.Range("K100").Value = 1566
$ws.Range("L100").Value = 4999
$ws.Range("M100").Value = -1025
$ws.Range("N100").Value = -6081

$ws.Range("H132").Value = 69317.07000000001
$ws.Range("I132").Value = 114296.89
$ws.Range("J132").Value = 1847.3334
$ws.Range("K132").Value = 342890.67
$ws.Range("L132").Value = 5542.0002
$ws.Range("M132").Value = -340360.67
$ws.Range("N132").Value = -10602.0002

$ws.Range("H133").Value = 89994.5
$ws.Range("J133").Value = 89994.5
$ws.Range("L133").Value = 89994.5
$ws.Range("N133").Value = -100114.5
